# RAMP workbook update: replace the "Households" user entry with a
# "Church" entry, and add a new "Ch_indoor_bulb" appliance row for it.

$wb = $excel.ActiveWorkbook

# --- Sheet "Users": rename Households -> Church, update N.users / preference
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Church"
$wsUsers.Range("B2").Value = 3
$wsUsers.Range("C2").Value = 0

# --- Sheet "Appliances": new appliance row (row 5) for the Church user
$wsApp = $wb.Worksheets.Item("Appliances")
$wsApp.Range("A5").Value = "Ch_indoor_bulb"
$wsApp.Range("B5").Value = "Church"
$wsApp.Range("C5").Value = 26
$wsApp.Range("E5").Value = 10
$wsApp.Range("F5").Value = 1
$wsApp.Range("G5").Value = 210
$wsApp.Range("H5").Value = 0.2
$wsApp.Range("I5").Value = 60
$wsApp.Range("J5").Value = "yes"
$wsApp.Range("M5").Value = "yes"
$wsApp.Range("P5").Value = "1200,1400"
$wsApp.Range("Q5").Value = "0,0"
$wsApp.Range("S5").Value = 0.1

# --- Selections / active sheet to match the saved view state
$wsUsers.Range("A2").Select() | Out-Null
$wsApp.Range("T5").Select() | Out-Null
$wsApp.Activate() | Out-Null
